$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns: Wins (AD), Losses (AE), Ties (AF) - the team's season record,
# appended to the right of the existing stats table.

# Header row (row 1): set the header text, then copy the formatting from the
# neighboring header cell (AC1) so the new headers match the bold/bordered
# header style used by the rest of row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows: every player row on this sheet shares the team's season record.
$wins = 96
$losses = 66
$ties = 0

$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
